$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.2305699481865285
$ws.Range("C2").Value = 0.4792746113989637
$ws.Range("J2").Value = 0.0155440414507772
$ws.Range("P2").Value = 0.1787564766839378
$ws.Range("S2").Value = 0.09585492227979274
$ws.Range("B3").Value = 0.02040816326530612
$ws.Range("C3").Value = 0.01530612244897959
$ws.Range("J3").Value = 0.02551020408163265
$ws.Range("P3").Value = 0.673469387755102
$ws.Range("S3").Value = 0.2653061224489796
$ws.Range("J4").Value = 0.0759493670886076
$ws.Range("O4").Value = 0.01265822784810127
$ws.Range("P4").Value = 0.6075949367088608
$ws.Range("S4").Value = 0.3037974683544304
$ws.Range("B6").Value = 0.09465020576131687
$ws.Range("D6").Value = 0.01646090534979424
$ws.Range("F6").Value = 0.09876543209876543
$ws.Range("J6").Value = 0.2345679012345679
$ws.Range("O6").Value = 0.02880658436213992
$ws.Range("Q6").Value = 0.1769547325102881
$ws.Range("R6").Value = 0.0411522633744856
$ws.Range("S6").Value = 0.308641975308642
$ws.Range("B7").Value = 0.1162790697674419
$ws.Range("D7").Value = 0.01937984496124031
$ws.Range("F7").Value = 0.06201550387596899
$ws.Range("J7").Value = 0.124031007751938
$ws.Range("O7").Value = 0.02713178294573643
$ws.Range("Q7").Value = 0.1744186046511628
$ws.Range("R7").Value = 0.1085271317829457
$ws.Range("S7").Value = 0.3682170542635659
$ws.Range("B8").Value = 0.1045364891518738
$ws.Range("D8").Value = 0.02761341222879684
$ws.Range("E8").Value = 0.001972386587771203
$ws.Range("F8").Value = 0.0670611439842209
$ws.Range("J8").Value = 0.1045364891518738
$ws.Range("O8").Value = 0.02761341222879684
$ws.Range("Q8").Value = 0.1755424063116371
$ws.Range("R8").Value = 0.09072978303747535
$ws.Range("S8").Value = 0.4003944773175542
$ws.Range("B9").Value = 0.1099476439790576
$ws.Range("D9").Value = 0.02617801047120419
$ws.Range("F9").Value = 0.07329842931937172
$ws.Range("J9").Value = 0.1047120418848168
$ws.Range("O9").Value = 0.01570680628272251
$ws.Range("Q9").Value = 0.225130890052356
$ws.Range("R9").Value = 0.06806282722513089
$ws.Range("S9").Value = 0.3769633507853403
$ws.Range("B10").Value = 0.1180705543556515
$ws.Range("D10").Value = 0.04031677465802735
$ws.Range("E10").Value = 0.0007199424046076314
$ws.Range("F10").Value = 0.06047516198704104
$ws.Range("J10").Value = 0.1065514758819294
$ws.Range("O10").Value = 0.01439884809215263
$ws.Range("Q10").Value = 0.2397408207343412
$ws.Range("R10").Value = 0.06911447084233262
$ws.Range("S10").Value = 0.3506119510439165
$ws.Range("G11").Value = 0.1089588377723971
$ws.Range("J11").Value = 0.1186440677966102
$ws.Range("K11").Value = 0.1864406779661017
$ws.Range("L11").Value = 0.5738498789346247
$ws.Range("S11").Value = 0.01210653753026634
$ws.Range("G12").Value = 0.744
$ws.Range("J12").Value = 0.184
$ws.Range("K12").Value = 0.012
$ws.Range("L12").Value = 0.036
$ws.Range("S12").Value = 0.024
$ws.Range("G13").Value = 0.7826086956521739
$ws.Range("J13").Value = 0.1739130434782609
$ws.Range("S13").Value = 0.04347826086956522
$ws.Range("F15").Value = 0.00423728813559322
$ws.Range("H15").Value = 0.1440677966101695
$ws.Range("I15").Value = 0.05508474576271186
$ws.Range("J15").Value = 0.3644067796610169
$ws.Range("K15").Value = 0.04661016949152542
$ws.Range("M15").Value = 0.01694915254237288
$ws.Range("O15").Value = 0.04661016949152542
$ws.Range("S15").Value = 0.3220338983050847
$ws.Range("F16").Value = 0.01239669421487603
$ws.Range("H16").Value = 0.2148760330578512
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.3677685950413223
$ws.Range("K16").Value = 0.1033057851239669
$ws.Range("M16").Value = 0.02892561983471074
$ws.Range("N16").Value = 0.004132231404958678
$ws.Range("O16").Value = 0.04132231404958678
$ws.Range("S16").Value = 0.1363636363636364
$ws.Range("F17").Value = 0.02550091074681239
$ws.Range("H17").Value = 0.1876138433515483
$ws.Range("I17").Value = 0.05828779599271403
$ws.Range("J17").Value = 0.4116575591985428
$ws.Range("K17").Value = 0.1238615664845173
$ws.Range("M17").Value = 0.01821493624772313
$ws.Range("O17").Value = 0.07650273224043716
$ws.Range("S17").Value = 0.09836065573770492
$ws.Range("F18").Value = 0.02116402116402116
$ws.Range("H18").Value = 0.1534391534391534
$ws.Range("I18").Value = 0.0582010582010582
$ws.Range("J18").Value = 0.4708994708994709
$ws.Range("K18").Value = 0.164021164021164
$ws.Range("M18").Value = 0.01058201058201058
$ws.Range("N18").Value = 0.005291005291005291
$ws.Range("O18").Value = 0.04232804232804233
$ws.Range("S18").Value = 0.07407407407407407
$ws.Range("F19").Value = 0.01404286770140429
$ws.Range("H19").Value = 0.2121212121212121
$ws.Range("I19").Value = 0.08795269770879527
$ws.Range("J19").Value = 0.3562453806356246
$ws.Range("K19").Value = 0.1433850702143385
$ws.Range("M19").Value = 0.02069475240206948
$ws.Range("N19").Value = 0.00147819660014782
$ws.Range("O19").Value = 0.06134515890613452
$ws.Range("S19").Value = 0.1027346637102735
